$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds dates formatted as plain text (DD/MM/YYYY strings), not
# real Excel dates. Force text number-format before assigning so COM
# doesn't auto-convert the literal into a date serial, then restore the
# default "Normal" style so no stray number-format style lingers on the
# cell (matches the rest of the column, which carries no explicit style).
$ws.Range("B74:B75").NumberFormat = "@"

# Row 74 (A74=73): 10/06/2020
$ws.Range("A74").Value = 73
$ws.Range("B74").Value = "10/06/2020"
$ws.Range("C74").Value = 237
$ws.Range("D74").Value = 6
$ws.Range("E74").Value = 46
$ws.Range("F74").Value = "156,0360265"
$ws.Range("G74").Value = "0,0253164557"
$ws.Range("H74").Value = 184
$ws.Range("I74").Value = 410
$ws.Range("J74").Value = 647
$ws.Range("K74").Value = 10
$ws.Range("L74").Value = 50
$ws.Range("M74").Value = 5
$ws.Range("N74").Value = 45
$ws.Range("O74").Value = 23
$ws.Range("P74").Value = 19
$ws.Range("Q74").Value = 9
$ws.Range("R74").Value = 11

# Row 75 (A75=74): 11/06/2020
$ws.Range("A75").Value = 74
$ws.Range("B75").Value = "11/06/2020"
$ws.Range("C75").Value = 246
$ws.Range("D75").Value = 7
$ws.Range("E75").Value = 42
$ws.Range("F75").Value = "161,9614453"
$ws.Range("G75").Value = "0,02845528455"
$ws.Range("H75").Value = 196
$ws.Range("I75").Value = 434
$ws.Range("J75").Value = 680
$ws.Range("K75").Value = 9
$ws.Range("L75").Value = 48
$ws.Range("M75").Value = 5
$ws.Range("N75").Value = 43
$ws.Range("O75").Value = 33
$ws.Range("P75").Value = 21
$ws.Range("Q75").Value = 9
$ws.Range("R75").Value = 11

$ws.Range("B74:B75").Style = "Normal"
